$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddProduct")

# Update the VAT values (column M) for rows 2-4
$ws.Range("M2").Value = 5000
$ws.Range("M3").Value = 2000
$ws.Range("M4").Value = 3000

# Update the active selection on the sheet
$ws.Range("L5").Select()
